$wb = $excel.ActiveWorkbook

# "Set Values Here" sheet holds the weighting matrix used by the GRA-* sheets.
# Row 8 corresponds to "carbon tax revenue". Update Deficit Spending (C8) and
# Payroll Taxes (E8) weights from 0 to 5 (Household Taxes / D8 stays at 5).
$setValues = $wb.Worksheets.Item("Set Values Here")
$setValues.Range("C8").Value = 5
$setValues.Range("E8").Value = 5

# Leave the selection on that sheet where it was left after editing the value (C9).
$setValues.Activate()
$setValues.Range("C9").Select()

# The carbon tax sheet pulls these weights via a TRANSPOSE array formula, so
# its values recompute automatically. Just restore the last-used selection.
$carbonTax = $wb.Worksheets.Item("GRA-carbontax")
$carbonTax.Activate()
$carbonTax.Range("B5").Select()

# Restore the originally active sheet/tab.
$wb.Worksheets.Item("About").Activate()
